# Add new row under range 301-450 in Card18
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

# Insert a new row at row 5, pushing existing rows 5-14 down to 6-15
$ws.Rows.Item(5).Insert()

# Fill in the new row 5 with the service-card range info.
# Force text formatting on A:C so the numbers are stored as text, matching
# the rest of the sheet (every other cell in this table is text, not numeric).
$ws.Range("A5:C5").NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "18"
$ws.Cells.Item(5, 2).Value = "301"
$ws.Cells.Item(5, 3).Value = "450"
$ws.Cells.Item(5, 15).Value = "محمد عبدالله"

# The row that was previously row 14 (now row 15) had its D:N cells blank;
# after the insert they should show the placeholder text "nan" like the
# rest of the sheet's unfilled cells.
$ws.Range("D15:N15").Value = "nan"

$wb.Save()
